$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.509.53'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.11'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '262.45'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5329'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3107'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -4.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06892'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.59'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7628'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07828'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.871.44'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.75'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.043'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.04'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007946'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.535.09'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.079.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.632'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.015'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.319'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.86'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.186'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.05'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.41'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.286'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08798'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.097'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04834'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.936'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7334'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.135'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.107'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.334'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +6.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01726'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4815'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9028'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '108.43'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.903'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.81%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.647'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4164'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.039'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1240'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.99'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.8999'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05805'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.97%  '
